$wb = $excel.ActiveWorkbook

# "se hace el ajuste de los perfiles para agendar asesoria quitando el
#  super chanchullo que hizo nuestro amigo sebastian"
# Add a new "asesoria" (advisory session) record to the asesorias sheet.
$ws = $wb.Worksheets.Item("asesorias")

$ws.Range("A3").Value = "Sebastian Palacio"
$ws.Range("B3").Value = "Juan Carlos Gil"
$ws.Range("C3").Value = "Consulta sobre módulos"

# "07-11-2023" parses as a valid date (dd-mm-yyyy) and Excel would silently
# convert it to a serial date, unlike the sibling "31-10-2023" cell above it
# (day 31 isn't a valid month, so it is stored as plain text). Round-trip it
# through a text formula so it lands back in the cell as a literal string,
# matching the existing date column's storage, without stamping a new
# NumberFormat style onto the cell.
$ws.Range("D3").Formula = '="07-11-2023"'
$ws.Range("D3").Value = $ws.Range("D3").Value

$ws.Range("E3").Value = "02:00 - 02:20"

# Move the active selection/tab over to "estudiantes" before saving.
$ws2 = $wb.Worksheets.Item("estudiantes")
$ws2.Activate()
$ws2.Range("A1").Select()
